$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New numeric header cells (row 1) ---
$ws.Range("I1").Value = 45
$ws.Range("K1").Value = 46

# --- New TO DO / FEEDBACK block + new text columns (authoring order) ---
$ws.Range("A10").Value = "TO DO/UPDATE/FEEDBACK"

$ws.Range("I4").Value = "Prepare Slide for Smart Home"

$ws.Range("I10").Value = "Need to add advantages & disadvatages"
$ws.Range("I11").Value = "Restructure it become more abstract"
$ws.Range("I12").Value = "Make points & reduce sentences"
$ws.Range("I13").Value = "Add more references"

$ws.Range("K5").Value = "Add disadvantage "
$ws.Range("K6").Value = "Add and edit references"
$ws.Range("K7").Value = "Add Advantages "
$ws.Range("K4").Value = "Restructure the slide"

$ws.Range("I5").Value = "Prepare Slide for Smart Home"
$ws.Range("I6").Value = "Prepare Slide for Smart Home"
$ws.Range("I7").Value = "Prepare Slide for Smart Home"

# --- Apply the built-in "Good" cell style to the data table cells ---
foreach ($addr in @("C1","E1","G1","I1","K1")) {
    $ws.Range($addr).Style = "Good"
}
foreach ($addr in @("C4","E4","G4","I4","K4")) {
    $ws.Range($addr).Style = "Good"
}
foreach ($addr in @("C5","E5","G5","I5","K5")) {
    $ws.Range($addr).Style = "Good"
}
foreach ($addr in @("C6","E6","G6","I6","K6")) {
    $ws.Range($addr).Style = "Good"
}
foreach ($addr in @("C7","E7","G7","I7","K7")) {
    $ws.Range($addr).Style = "Good"
}

# --- Apply the built-in "Neutral" cell style to the TO DO / FEEDBACK block ---
foreach ($addr in @("A10","C10","E10","G10","I10","I11","I12","I13")) {
    $ws.Range($addr).Style = "Neutral"
}

# --- Column widths for the new columns ---
$ws.Range("A1").EntireColumn.ColumnWidth = 36.83
$ws.Range("I1").EntireColumn.ColumnWidth = 32.94
$ws.Range("K1").EntireColumn.ColumnWidth = 24.72

# --- View: select H10 like the saved workbook ---
$ws.Range("H10").Select()
